# Assignment1.xlsx update — "Add files via upload"
# Applies the data edits made to the "PO List" sheet plus the
# conditional-formatting rule tweak on F3:F29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO List")

# ---------------------------------------------------------------
# Cell value updates (column N "10storey_next" refresh, plus a
# handful of related columns on rows 7 and 20).
# ---------------------------------------------------------------
$ws.Range("N3").Value = 13
$ws.Range("N4").Value = 20
$ws.Range("N5").Value = 23
$ws.Range("N6").Value = 17

$ws.Range("K7").Value = 5
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = 44923
$ws.Range("N7").Value = 1

$ws.Range("N8").Value = 22
$ws.Range("N9").Value = 6
$ws.Range("N10").Value = 12
$ws.Range("N11").Value = 10
$ws.Range("N12").Value = 8
$ws.Range("N13").Value = 7
$ws.Range("N14").Value = 21
$ws.Range("N15").Value = 14
$ws.Range("N16").Value = 16
$ws.Range("N17").Value = 15
$ws.Range("N18").Value = 24
$ws.Range("N19").Value = 4

$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 2
$ws.Range("M20").Value = 44923
$ws.Range("N20").Value = 1
$ws.Range("S20").Value = 2
$ws.Range("T20").Value = 2
$ws.Range("U20").Value = 44923

$ws.Range("N21").Value = 19
$ws.Range("N24").Value = 9
$ws.Range("N25").Value = 5
$ws.Range("N26").Value = 17
$ws.Range("N27").Value = 3
$ws.Range("N29").Value = 11

# ---------------------------------------------------------------
# Active-cell selection moved from P40 to E36 in the bottom-right
# (frozen) pane of the "PO List" sheet view.
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("E36").Select()

# ---------------------------------------------------------------
# Conditional formatting on F3:F29: re-create the "less than 3"
# highlight rule so a second (duplicate) dxf is registered and the
# rule now points at it, keeping the same red-on-red look while
# bumping dxfs count from 1 to 2.
# ---------------------------------------------------------------
$cfRange = $ws.Range("F3:F29")
$cfRange.FormatConditions.Delete()
$newRule = $cfRange.FormatConditions.Add(1, 6, "3")
$newRule.Font.Color = 393372
$newRule.Interior.Color = 13551615
$newRule.Priority = 6
